# agregaProducto: append a new product row right after the last used row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$producto = @{
    Id       = "12314"
    Nombre   = "prueba"
    Precio   = "12313"
    Cantidad = "12313"
}

$ws.Cells.Item($newRow, 1).Value = $producto.Id
$ws.Cells.Item($newRow, 2).Value = $producto.Nombre
$ws.Cells.Item($newRow, 3).Value = $producto.Precio
$ws.Cells.Item($newRow, 4).Value = $producto.Cantidad

# Keep the new row's formatting consistent with the other plain data rows
# (no explicit style / text-only column format), matching how the rest of
# the "clientes" data rows already look.
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4)).Style = "Normal"
